$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cryptocurrency price / volume figures scraped on
# Fri Mar  1 22:38:56 UTC 2024.
$ws.Range("D2").Value = "'62.406.50"
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = "'3.431.92"
$ws.Range("E3").Value = '  +3.35%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'407.39"
$ws.Range("E5").Value = '  +2.55%  '
$ws.Range("D6").Value = "'130.24"
$ws.Range("E6").Value = '  +4.68%  '
$ws.Range("D7").Value = "'0.601"
$ws.Range("E7").Value = '  +3.10%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = "'0.697"
$ws.Range("E9").Value = '  +7.29%  '
$ws.Range("D10").Value = "'0.144"
$ws.Range("E10").Value = '  +23.06%  '
$ws.Range("D11").Value = "'42.17"
$ws.Range("E11").Value = '  +4.48%  '
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = "'8.54"
$ws.Range("E13").Value = '  +5.01%  '
$ws.Range("E14").Value = '  +4.96%  '
$ws.Range("D15").Value = "'3.421.01"
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = "'62.474.09"
$ws.Range("E16").Value = '  +2.20%  '
$ws.Range("D17").Value = "'11.57"
$ws.Range("E17").Value = '  +4.92%  '
$ws.Range("D18").Value = "'0.0000168"
$ws.Range("E18").Value = '  +35.98%  '
$ws.Range("E19").Value = '  +3.11%  '
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("D21").Value = "'84.86"
$ws.Range("E21").Value = '  +7.48%  '
$ws.Range("D22").Value = "'314.75"
$ws.Range("E22").Value = '  +6.58%  '
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("E24").Value = '  +4.09%  '
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").Value = "'29.96"
$ws.Range("E26").Value = '  +4.89%  '
$ws.Range("E27").Value = '  +1.69%  '
$ws.Range("D28").Value = "'7.78"
$ws.Range("E28").Value = '  +5.97%  '
$ws.Range("D29").Value = "'2.71"
$ws.Range("E29").Value = '  +8.96%  '
$ws.Range("E30").Value = '  +2.06%  '
$ws.Range("D31").Value = "'44.22"
$ws.Range("E31").Value = '  +9.54%  '
$ws.Range("E32").Value = '  +2.99%  '
$ws.Range("E33").Value = '  +3.15%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = "'0.0486"
$ws.Range("E35").Value = '  +3.49%  '
$ws.Range("D36").Value = "'51.52"
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("E38").Value = '  +4.04%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = "'3.33"
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = "'0.319"
$ws.Range("E40").Value = '  +16.21%  '
$ws.Range("D41").Value = "'143.94"
$ws.Range("E41").Value = '  +5.93%  '
$ws.Range("E42").Value = '  +4.63%  '
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("D44").Value = "'16.94"
$ws.Range("E44").Value = '  +3.45%  '
$ws.Range("D45").Value = "'3.93"
$ws.Range("E45").Value = '  +3.83%  '
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").Value = "'21.35"
$ws.Range("E47").Value = '  +2.29%  '
$ws.Range("D48").Value = "'2.108.31"
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("E49").Value = '  +11.88%  '
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").Value = "'1.08"
$ws.Range("E51").Value = '  +32.05%  '

